# Update column F ("dSF") values for specific rows to match repulled data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = -2
    3  = -10
    4  = 1
    7  = 0
    11 = -2
    16 = -4
    20 = -1
    22 = -4
    25 = -8
    31 = 6
    33 = -3
    34 = -5
    38 = -2
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
